$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.157854437828064
$ws.Range("B1").Value = 2.384793281555176
$ws.Range("D1").Value = 2.388535737991333
$ws.Range("E1").Value = 1.222447395324707
